$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.417859666666667
$ws.Range("H2").Value = 4.253579
$ws.Range("I2").Value = 0.1472651073415806
$ws.Range("J2").Value = 0.1472651073415806
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.746298666666667
$ws.Range("N2").Value = 11.238896
$ws.Range("O2").Value = 0.2944715732214049
$ws.Range("P2").Value = 0.294471573221405
$ws.Range("Q2").Value = 5.311725778753778
$ws.Range("R2").Value = 47.80553200878401
$ws.Range("S2").Value = 0.04336538783949431
$ws.Range("T2").Value = 0.04336538783949433

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.417859666666667
$ws.Range("H3").Value = 4.253579
$ws.Range("I3").Value = 0.1472651073415806
$ws.Range("J3").Value = 0.1472651073415806
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.370261666666667
$ws.Range("N3").Value = 13.110785
$ws.Range("O3").Value = 0.3435171466234404
$ws.Range("P3").Value = 0.3435171466234404
$ws.Range("Q3").Value = 6.196417749946112
$ws.Range("R3").Value = 55.767759749515
$ws.Range("S3").Value = 0.05058808947117444
$ws.Range("T3").Value = 0.05058808947117444

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.417859666666667
$ws.Range("H4").Value = 4.253579
$ws.Range("I4").Value = 0.1472651073415806
$ws.Range("J4").Value = 0.1472651073415806
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.837969333333333
$ws.Range("N4").Value = 5.513908
$ws.Range("O4").Value = 0.1444705212467569
$ws.Range("P4").Value = 0.1444705212467569
$ws.Range("Q4").Value = 2.605982586303556
$ws.Range("R4").Value = 23.453843276732
$ws.Range("S4").Value = 0.02127546681909775
$ws.Range("T4").Value = 0.02127546681909775

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.417859666666667
$ws.Range("H5").Value = 4.253579
$ws.Range("I5").Value = 0.1472651073415806
$ws.Range("J5").Value = 0.1472651073415806
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.767576666666667
$ws.Range("N5").Value = 8.30273
$ws.Range("O5").Value = 0.2175407589083977
$ws.Range("P5").Value = 0.2175407589083977
$ws.Range("Q5").Value = 3.924035330074445
$ws.Range("R5").Value = 35.31631797067001
$ws.Range("S5").Value = 0.0320361632118141
$ws.Range("T5").Value = 0.03203616321181411

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.189892666666666
$ws.Range("H6").Value = 18.569678
$ws.Range("I6").Value = 0.6429093297593833
$ws.Range("J6").Value = 0.6429093297593833
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.746298666666667
$ws.Range("N6").Value = 11.238896
$ws.Range("O6").Value = 0.2944715732214049
$ws.Range("P6").Value = 0.294471573221405
$ws.Range("Q6").Value = 23.18918664394311
$ws.Range("R6").Value = 208.702679795488
$ws.Range("S6").Value = 0.1893185217729646
$ws.Range("T6").Value = 0.1893185217729647

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.189892666666666
$ws.Range("H7").Value = 18.569678
$ws.Range("I7").Value = 0.6429093297593833
$ws.Range("J7").Value = 0.6429093297593833
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.370261666666667
$ws.Range("N7").Value = 13.110785
$ws.Range("O7").Value = 0.3435171466234404
$ws.Range("P7").Value = 0.3435171466234404
$ws.Range("Q7").Value = 27.05145064191445
$ws.Range("R7").Value = 243.46305577723
$ws.Range("S7").Value = 0.2208503784965319
$ws.Range("T7").Value = 0.2208503784965319

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.189892666666666
$ws.Range("H8").Value = 18.569678
$ws.Range("I8").Value = 0.6429093297593833
$ws.Range("J8").Value = 0.6429093297593833
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.837969333333333
$ws.Range("N8").Value = 5.513908
$ws.Range("O8").Value = 0.1444705212467569
$ws.Range("P8").Value = 0.1444705212467569
$ws.Range("Q8").Value = 11.37683289795822
$ws.Range("R8").Value = 102.391496081624
$ws.Range("S8").Value = 0.0928814459847412
$ws.Range("T8").Value = 0.0928814459847412

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.189892666666666
$ws.Range("H9").Value = 18.569678
$ws.Range("I9").Value = 0.6429093297593833
$ws.Range("J9").Value = 0.6429093297593833
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.767576666666667
$ws.Range("N9").Value = 8.30273
$ws.Range("O9").Value = 0.2175407589083977
$ws.Range("P9").Value = 0.2175407589083977
$ws.Range("Q9").Value = 17.13100251343778
$ws.Range("R9").Value = 154.17902262094
$ws.Range("S9").Value = 0.1398589835051456
$ws.Range("T9").Value = 0.1398589835051456

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.012114666666667
$ws.Range("H10").Value = 3.036344
$ws.Range("I10").Value = 0.105122656728831
$ws.Range("J10").Value = 0.105122656728831
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.746298666666667
$ws.Range("N10").Value = 11.238896
$ws.Range("O10").Value = 0.2944715732214049
$ws.Range("P10").Value = 0.294471573221405
$ws.Range("Q10").Value = 3.791683826247111
$ws.Range("R10").Value = 34.125154436224
$ws.Range("S10").Value = 0.03095563410815258
$ws.Range("T10").Value = 0.03095563410815258

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.012114666666667
$ws.Range("H11").Value = 3.036344
$ws.Range("I11").Value = 0.105122656728831
$ws.Range("J11").Value = 0.105122656728831
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.370261666666667
$ws.Range("N11").Value = 13.110785
$ws.Range("O11").Value = 0.3435171466234404
$ws.Range("P11").Value = 0.3435171466234404
$ws.Range("Q11").Value = 4.423205930004445
$ws.Range("R11").Value = 39.80885337004
$ws.Range("S11").Value = 0.03611143508496344
$ws.Range("T11").Value = 0.03611143508496344

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.012114666666667
$ws.Range("H12").Value = 3.036344
$ws.Range("I12").Value = 0.105122656728831
$ws.Range("J12").Value = 0.105122656728831
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.837969333333333
$ws.Range("N12").Value = 5.513908
$ws.Range("O12").Value = 0.1444705212467569
$ws.Range("P12").Value = 0.1444705212467569
$ws.Range("Q12").Value = 1.860235719150222
$ws.Range("R12").Value = 16.742121472352
$ws.Range("S12").Value = 0.01518712501245811
$ws.Range("T12").Value = 0.01518712501245811

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.012114666666667
$ws.Range("H13").Value = 3.036344
$ws.Range("I13").Value = 0.105122656728831
$ws.Range("J13").Value = 0.105122656728831
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.767576666666667
$ws.Range("N13").Value = 8.30273
$ws.Range("O13").Value = 0.2175407589083977
$ws.Range("P13").Value = 0.2175407589083977
$ws.Range("Q13").Value = 2.801104935457778
$ws.Range("R13").Value = 25.20994441912
$ws.Range("S13").Value = 0.02286846252325688
$ws.Range("T13").Value = 0.02286846252325688

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.008073333333333
$ws.Range("H14").Value = 3.02422
$ws.Range("I14").Value = 0.1047029061702051
$ws.Range("J14").Value = 0.1047029061702051
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.746298666666667
$ws.Range("N14").Value = 11.238896
$ws.Range("O14").Value = 0.2944715732214049
$ws.Range("P14").Value = 0.294471573221405
$ws.Range("Q14").Value = 3.776543784568889
$ws.Range("R14").Value = 33.98889406112
$ws.Range("S14").Value = 0.03083202950079345
$ws.Range("T14").Value = 0.03083202950079346

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.008073333333333
$ws.Range("H15").Value = 3.02422
$ws.Range("I15").Value = 0.1047029061702051
$ws.Range("J15").Value = 0.1047029061702051
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.370261666666667
$ws.Range("N15").Value = 13.110785
$ws.Range("O15").Value = 0.3435171466234404
$ws.Range("P15").Value = 0.3435171466234404
$ws.Range("Q15").Value = 4.405544245855556
$ws.Range("R15").Value = 39.6498982127
$ws.Range("S15").Value = 0.03596724357077068
$ws.Range("T15").Value = 0.03596724357077068

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.008073333333333
$ws.Range("H16").Value = 3.02422
$ws.Range("I16").Value = 0.1047029061702051
$ws.Range("J16").Value = 0.1047029061702051
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.837969333333333
$ws.Range("N16").Value = 5.513908
$ws.Range("O16").Value = 0.1444705212467569
$ws.Range("P16").Value = 0.1444705212467569
$ws.Range("Q16").Value = 1.852807872417778
$ws.Range("R16").Value = 16.67527085176
$ws.Range("S16").Value = 0.01512648343045981
$ws.Range("T16").Value = 0.01512648343045981

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.008073333333333
$ws.Range("H17").Value = 3.02422
$ws.Range("I17").Value = 0.1047029061702051
$ws.Range("J17").Value = 0.1047029061702051
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.767576666666667
$ws.Range("N17").Value = 8.30273
$ws.Range("O17").Value = 0.2175407589083977
$ws.Range("P17").Value = 0.2175407589083977
$ws.Range("Q17").Value = 2.789920235622223
$ws.Range("R17").Value = 25.1092821206
$ws.Range("S17").Value = 0.02277714966818118
$ws.Range("T17").Value = 0.02277714966818118
